$wb = $excel.ActiveWorkbook

# Sheet "arbolts": update B2:E2 (A2 stays 3)
$ws1 = $wb.Worksheets.Item("arbolts")
$ws1.Range("B2").Value = 11.74729444162726
$ws1.Range("C2").Value = 270.5656953227698
$ws1.Range("D2").Value = 16.44888127876087
$ws1.Range("E2").Value = 0.8726742175255467

# Sheet "bosquets": update A2:E2
$ws2 = $wb.Worksheets.Item("bosquets")
$ws2.Range("A2").Value = 102
$ws2.Range("B2").Value = 8.779559828778744
$ws2.Range("C2").Value = 142.8535896048408
$ws2.Range("D2").Value = 11.95213744921137
$ws2.Range("E2").Value = 0.9327743857031751

# Sheet "knnts": update A2:E2
$ws3 = $wb.Worksheets.Item("knnts")
$ws3.Range("A2").Value = 27
$ws3.Range("B2").Value = 9.068932712817771
$ws3.Range("C2").Value = 156.0077977761895
$ws3.Range("D2").Value = 12.49030815377225
$ws3.Range("E2").Value = 0.9265841336601331
